$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (ano/ano_obj = 2025) with refreshed metrics
$ws.Range("C8").Value = 1037
$ws.Range("D8").Value = 172
$ws.Range("E8").Value = 865
$ws.Range("F8").Value = 7.054963084495489
$ws.Range("G8").Value = 83.41369334619093
$ws.Range("H8").Value = 16.58630665380906
